# Tag Excel List v1.0 -> v1.1 : add a "detagging" function
# Adds three new tag rows for "Leck, Frances O." so alternate tag
# spellings ("렉 부인 2", "Mrs. 렉", "김 렉") resolve to the same person.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New tag -> full-name rows appended right after the existing data (row 40).
$ws.Range("A41").Value = "렉 부인 2"
$ws.Range("B41").Value = "Leck, Frances O."

$ws.Range("A42").Value = "Mrs. 렉"
$ws.Range("B42").Value = "Leck, Frances O."

$ws.Range("A43").Value = "김 렉"
$ws.Range("B43").Value = "Leck, Frances O."

# Column A keeps the same centered style used by the rest of the list.
$ws.Range("A41:A43").HorizontalAlignment = -4108
$ws.Range("A41:A43").VerticalAlignment = -4108

# Scroll the view down and leave the selection on the next empty row,
# matching where the editor's cursor ended up after the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 23
$ws.Range("A44").Select()
